$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 1991
$ws.Range("E2").Value = 96
$ws.Range("F2").Value = 96
$ws.Range("G2").Value = 139
$ws.Range("H2").Value = 102
$ws.Range("I2").Value = 102
$ws.Range("K2").Value = 2990
$ws.Range("L2").Value = 765
$ws.Range("M2").Value = 2225
$ws.Range("N2").Value = 2225
$ws.Range("P2").Value = 150
$ws.Range("Q2").Value = -23
$ws.Range("R2").Value = -35
$ws.Range("S2").Value = 75
$ws.Range("T2").Value = 51
$ws.Range("U2").Value = -74
$ws.Range("V2").Value = 489
$ws.Range("W2").Value = 4.8
$ws.Range("X2").Value = 5.13
$ws.Range("Y2").Value = 4.68
$ws.Range("Z2").Value = 3.57
$ws.Range("AA2").Value = 34.36
$ws.Range("AB2").Value = 1380.47
$ws.Range("AC2").Value = 3405
$ws.Range("AD2").Value = 25.38
$ws.Range("AE2").Value = 74181
$ws.Range("AF2").Value = 1.16
$ws.Range("AG2").Value = 750
$ws.Range("AH2").Value = 0.87
$ws.Range("AI2").Value = 22.03
$ws.Range("AJ2").Value = 3000000

# Row 3
$ws.Range("D3").Value = 2150
$ws.Range("E3").Value = 120
$ws.Range("F3").Value = 120
$ws.Range("G3").Value = 139
$ws.Range("H3").Value = 105
$ws.Range("I3").Value = 105
$ws.Range("K3").Value = 3137
$ws.Range("L3").Value = 831
$ws.Range("M3").Value = 2306
$ws.Range("N3").Value = 2306
$ws.Range("P3").Value = 150
$ws.Range("Q3").Value = -205
$ws.Range("R3").Value = 116
$ws.Range("S3").Value = 70
$ws.Range("T3").Value = 53
$ws.Range("U3").Value = -258
$ws.Range("V3").Value = 583
$ws.Range("W3").Value = 5.6
$ws.Range("X3").Value = 4.9
$ws.Range("Y3").Value = 4.65
$ws.Range("Z3").Value = 3.44
$ws.Range("AA3").Value = 36.05
$ws.Range("AB3").Value = 1433.14
$ws.Range("AC3").Value = 3509
$ws.Range("AD3").Value = 13.81
$ws.Range("AE3").Value = 76858
$ws.Range("AF3").Value = 0.63
$ws.Range("AG3").Value = 825
$ws.Range("AH3").Value = 1.7
$ws.Range("AI3").Value = 23.51
$ws.Range("AJ3").Value = 3000000

# Row 4
$ws.Range("D4").Value = 1958
$ws.Range("E4").Value = 311
$ws.Range("F4").Value = 311
$ws.Range("G4").Value = 305
$ws.Range("H4").Value = 229
$ws.Range("I4").Value = 229
$ws.Range("K4").Value = 2892
$ws.Range("L4").Value = 393
$ws.Range("M4").Value = 2499
$ws.Range("N4").Value = 2499
$ws.Range("P4").Value = 150
$ws.Range("Q4").Value = 435
$ws.Range("R4").Value = 83
$ws.Range("S4").Value = -451
$ws.Range("T4").Value = 22
$ws.Range("U4").Value = 414
$ws.Range("V4").Value = 156
$ws.Range("W4").Value = 15.89
$ws.Range("X4").Value = 11.67
$ws.Range("Y4").Value = 9.51
$ws.Range("Z4").Value = 7.58
$ws.Range("AA4").Value = 15.72
$ws.Range("AB4").Value = 1564.42
$ws.Range("AC4").Value = 7618
$ws.Range("AD4").Value = 6.98
$ws.Range("AE4").Value = 83299
$ws.Range("AF4").Value = 0.64
$ws.Range("AG4").Value = 750
$ws.Range("AH4").Value = 1.41
$ws.Range("AI4").Value = 9.85
$ws.Range("AJ4").Value = 3000000

# Row 5
$ws.Range("D5").Value = 1706
$ws.Range("E5").Value = 302
$ws.Range("F5").Value = 302
$ws.Range("G5").Value = 333
$ws.Range("H5").Value = 249
$ws.Range("I5").Value = 249
$ws.Range("K5").Value = 3108
$ws.Range("L5").Value = 392
$ws.Range("M5").Value = 2716
$ws.Range("N5").Value = 2716
$ws.Range("P5").Value = 150
$ws.Range("Q5").Value = 327
$ws.Range("R5").Value = -274
$ws.Range("S5").Value = -36
$ws.Range("T5").Value = 28
$ws.Range("U5").Value = 299
$ws.Range("V5").Value = 137
$ws.Range("W5").Value = 17.69
$ws.Range("X5").Value = 14.61
$ws.Range("Y5").Value = 9.56
$ws.Range("Z5").Value = 8.31
$ws.Range("AA5").Value = 14.43
$ws.Range("AB5").Value = 1712.51
$ws.Range("AC5").Value = 8310
$ws.Range("AD5").Value = 4.95
$ws.Range("AE5").Value = 90524
$ws.Range("AF5").Value = 0.45
$ws.Range("AG5").Value = 750
$ws.Range("AH5").Value = 1.82
$ws.Range("AI5").Value = 9.03
$ws.Range("AJ5").Value = 3000000

# Row 6
$ws.Range("D6").Value = 1832
$ws.Range("E6").Value = 321
$ws.Range("F6").Value = 321
$ws.Range("G6").Value = 340
$ws.Range("H6").Value = 253
$ws.Range("I6").Value = 253
$ws.Range("K6").Value = 3275
$ws.Range("L6").Value = 347
$ws.Range("M6").Value = 2928
$ws.Range("N6").Value = 2928
$ws.Range("P6").Value = 150
$ws.Range("Q6").Value = 427
$ws.Range("R6").Value = -353
$ws.Range("S6").Value = -58
$ws.Range("T6").Value = 57
$ws.Range("U6").Value = 370
$ws.Range("V6").Value = 118
$ws.Range("W6").Value = 17.54
$ws.Range("X6").Value = 13.81
$ws.Range("Y6").Value = 8.96
$ws.Range("Z6").Value = 7.93
$ws.Range("AA6").Value = 11.85
$ws.Range("AB6").Value = 1866.08
$ws.Range("AC6").Value = 8432
$ws.Range("AD6").Value = 5.06
$ws.Range("AE6").Value = 99000
$ws.Range("AF6").Value = 0.43
$ws.Range("AG6").Value = 750
$ws.Range("AH6").Value = 1.76
$ws.Range("AI6").Value = 8.77
$ws.Range("AJ6").Value = 3000000

# Clear J and O columns for rows 2-5 (removed in target)
$ws.Range("J2").ClearContents()
$ws.Range("O2").ClearContents()
$ws.Range("J3").ClearContents()
$ws.Range("O3").ClearContents()
$ws.Range("J4").ClearContents()
$ws.Range("O4").ClearContents()
$ws.Range("J5").ClearContents()
$ws.Range("O5").ClearContents()

# Clear rows 7-9 data columns D:AJ entirely, keep only A,B,C
$ws.Range("D7:AJ7").ClearContents()
$ws.Range("D8:AJ8").ClearContents()
$ws.Range("D9:AJ9").ClearContents()
